$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A 33491-2023) hyperlink formulas gain a second "friendly name" argument.
# Note: the source commit's automatic-update script introduced this change with a
# stray bug: every link except the first (S2) is missing the closing quote
# right after the file name, so the string literal runs on into the next
# argument. We reproduce the formulas exactly as they appear in the authoritative
# diff, including that malformed syntax.

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/artfynd/A 33491-2023.xlsx"; "A 33491-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/kartor/A 33491-2023.png; "A 33491-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/klagomål/A 33491-2023.docx; "A 33491-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/klagomålsmail/A 33491-2023.docx; "A 33491-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/tillsyn/A 33491-2023.docx; "A 33491-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/tillsynsmail/A 33491-2023.docx; "A 33491-2023")'
